$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" (R20's old greeting) is replaced with the new commit note.
$ws.Range("E8").Value = "GIT UPDATE"

# Move/confirm the active selection on the edited cell, as Excel would after a manual edit.
$ws.Range("E8").Select()
